# Update the "想去人数" (F column) figures for the 江西-漫展信息 workbook
# to the newly scraped counts, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for sheet "展览"
$updatesExhibition = @{
    3  = 415
    5  = 126
    10 = 1177
    11 = 1479
    17 = 58
    19 = 262
    21 = 309
    25 = 169
    28 = 213
    29 = 4043
    32 = 246
    33 = 1044
    34 = 122
    36 = 295
    38 = 156
}

# Row -> new F-column value for sheet "全部类型"
$updatesAllTypes = @{
    3  = 415
    5  = 126
    10 = 1177
    11 = 1479
    17 = 58
    19 = 262
    21 = 309
    25 = 169
    28 = 213
    29 = 4043
    32 = 246
    33 = 1044
    34 = 122
    36 = 296
    38 = 156
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $updatesExhibition[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAllTypes.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $updatesAllTypes[$row]
}
